# Scheduled data refresh: update Leve profit-calculation columns (H:N)
# across several sheets to reflect the latest market board averages.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 211.96
$ws.Range("I15").Value = 211.96
$ws.Range("K15").Value = 635.88
$ws.Range("M15").Value = -466.88
$ws.Range("H40").Value = 0
$ws.Range("J40").Value = 0
$ws.Range("L40").Value = 0
$ws.Range("N40").ClearContents()
$ws.Range("H69").Value = 3500
$ws.Range("I69").Value = 3500
$ws.Range("K69").Value = 10500
$ws.Range("M69").Value = -9626
$ws.Range("H72").Value = 3500
$ws.Range("I72").Value = 3500
$ws.Range("K72").Value = 31500
$ws.Range("M72").Value = -27132
$ws.Range("H88").Value = 5199.4
$ws.Range("I88").Value = 998.5
$ws.Range("J88").Value = 8000
$ws.Range("K88").Value = 998.5
$ws.Range("L88").Value = 8000
$ws.Range("M88").Value = -592.5
$ws.Range("N88").Value = -8812
$ws.Range("H91").Value = 5199.4
$ws.Range("I91").Value = 998.5
$ws.Range("J91").Value = 8000
$ws.Range("K91").Value = 998.5
$ws.Range("L91").Value = 8000
$ws.Range("M91").Value = 405.5
$ws.Range("N91").Value = -10808
$ws.Range("H113").Value = 2268
$ws.Range("H135").Value = 683
$ws.Range("I135").Value = 683
$ws.Range("K135").Value = 6147
$ws.Range("M135").Value = -3612
$ws.Range("H137").Value = 4795.6924
$ws.Range("I137").Value = 4546.364
$ws.Range("K137").Value = 13639.092
$ws.Range("M137").Value = -11089.092

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H3").Value = 0
$ws.Range("J3").Value = 0
$ws.Range("L3").Value = 0
$ws.Range("N3").ClearContents()
$ws.Range("H13").Value = 2500250
$ws.Range("I13").Value = 5000000
$ws.Range("J13").Value = 500
$ws.Range("K13").Value = 5000000
$ws.Range("L13").Value = 500
$ws.Range("M13").Value = -4999856
$ws.Range("N13").Value = -788
$ws.Range("H32").Value = 117.666664
$ws.Range("I32").Value = 127
$ws.Range("J32").Value = 99
$ws.Range("K32").Value = 127
$ws.Range("L32").Value = 99
$ws.Range("M32").Value = 160
$ws.Range("N32").Value = -673
$ws.Range("H61").Value = 1759.4
$ws.Range("I61").Value = 1949.25
$ws.Range("K61").Value = 1949.25
$ws.Range("M61").Value = -1737.25
$ws.Range("H63").Value = 8237.25
$ws.Range("I63").Value = 8483.333
$ws.Range("J63").Value = 7499
$ws.Range("K63").Value = 8483.333
$ws.Range("L63").Value = 7499
$ws.Range("M63").Value = -7797.333000000001
$ws.Range("N63").Value = -8871
$ws.Range("H66").Value = 8237.25
$ws.Range("I66").Value = 8483.333
$ws.Range("J66").Value = 7499
$ws.Range("K66").Value = 42416.665
$ws.Range("L66").Value = 37495
$ws.Range("M66").Value = -38984.665
$ws.Range("N66").Value = -44359
$ws.Range("H74").Value = 3395
$ws.Range("I74").Value = 3193.3333
$ws.Range("K74").Value = 3193.3333
$ws.Range("M74").Value = -2319.3333
$ws.Range("H77").Value = 3395
$ws.Range("I77").Value = 3193.3333
$ws.Range("K77").Value = 15966.6665
$ws.Range("M77").Value = -11598.6665
$ws.Range("H97").Value = 949.6667
$ws.Range("I97").Value = 949.6667
$ws.Range("K97").Value = 949.6667
$ws.Range("M97").Value = -453.6667
$ws.Range("H110").Value = 550
$ws.Range("I110").Value = 1000
$ws.Range("J110").Value = 400
$ws.Range("K110").Value = 1000
$ws.Range("L110").Value = 400
$ws.Range("M110").Value = 1045
$ws.Range("N110").Value = -4490
$ws.Range("H136").Value = 1759.4
$ws.Range("I136").Value = 1949.25
$ws.Range("K136").Value = 5847.75
$ws.Range("M136").Value = -3297.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 297.33334
$ws.Range("I22").Value = 297.33334
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 297.33334
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = 52.66665999999998
$ws.Range("N22").ClearContents()
$ws.Range("H99").Value = 2000
$ws.Range("I99").Value = 2000
$ws.Range("K99").Value = 2000
$ws.Range("M99").Value = -502
$ws.Range("H126").Value = 2000
$ws.Range("I126").Value = 2000
$ws.Range("K126").Value = 6000
$ws.Range("M126").Value = -3530

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 61.18182
$ws.Range("I2").Value = 64.77778
$ws.Range("J2").Value = 45
$ws.Range("K2").Value = 388.66668
$ws.Range("L2").Value = 270
$ws.Range("M2").Value = -275.66668
$ws.Range("N2").Value = -496
$ws.Range("H8").Value = 925.25
$ws.Range("I8").Value = 925.25
$ws.Range("K8").Value = 2775.75
$ws.Range("M8").Value = -2636.75
$ws.Range("H80").Value = 1499.5
$ws.Range("J80").Value = 1499.5
$ws.Range("L80").Value = 4498.5
$ws.Range("N80").Value = -6370.5
$ws.Range("H83").Value = 1499.5
$ws.Range("J83").Value = 1499.5
$ws.Range("L83").Value = 13495.5
$ws.Range("N83").Value = -22855.5
$ws.Range("H107").Value = 738.73334
$ws.Range("I107").Value = 699.375
$ws.Range("K107").Value = 2098.125
$ws.Range("M107").Value = -178.125

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H45").Value = 70000
$ws.Range("J45").Value = 70000
$ws.Range("L45").Value = 70000
$ws.Range("N45").Value = -71118
$ws.Range("H109").Value = 25750
$ws.Range("I109").Value = 12500
$ws.Range("K109").Value = 12500
$ws.Range("M109").Value = -11460

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H12").Value = 3740.7144
$ws.Range("I12").Value = 4640
$ws.Range("J12").Value = 1492.5
$ws.Range("K12").Value = 4640
$ws.Range("L12").Value = 1492.5
$ws.Range("M12").Value = -4470
$ws.Range("N12").Value = -1832.5
$ws.Range("H40").Value = 0
$ws.Range("I40").Value = 0
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 0
$ws.Range("L40").Value = 0
$ws.Range("M40").ClearContents()
$ws.Range("N40").ClearContents()
$ws.Range("H46").Value = 5049.5
$ws.Range("I46").Value = 3862.6365
$ws.Range("J46").Value = 6053.769
$ws.Range("K46").Value = 3862.6365
$ws.Range("L46").Value = 6053.769
$ws.Range("M46").Value = -3674.6365
$ws.Range("N46").Value = -6429.769
$ws.Range("H68").Value = 2985.5
$ws.Range("I68").Value = 2647.3333
$ws.Range("J68").Value = 4000
$ws.Range("K68").Value = 2647.3333
$ws.Range("L68").Value = 4000
$ws.Range("M68").Value = -1898.3333
$ws.Range("N68").Value = -5498
$ws.Range("H71").Value = 2985.5
$ws.Range("I71").Value = 2647.3333
$ws.Range("J71").Value = 4000
$ws.Range("K71").Value = 13236.6665
$ws.Range("L71").Value = 20000
$ws.Range("M71").Value = -9492.6665
$ws.Range("N71").Value = -27488
$ws.Range("H82").Value = 1000
$ws.Range("I82").Value = 1000
$ws.Range("J82").Value = 0
$ws.Range("K82").Value = 1000
$ws.Range("L82").Value = 0
$ws.Range("M82").Value = -639
$ws.Range("N82").ClearContents()
$ws.Range("H85").Value = 1000
$ws.Range("I85").Value = 1000
$ws.Range("J85").Value = 0
$ws.Range("K85").Value = 1000
$ws.Range("L85").Value = 0
$ws.Range("M85").Value = 248
$ws.Range("N85").ClearContents()
$ws.Range("H130").Value = 69999
$ws.Range("J130").Value = 69999
$ws.Range("N130").Value = -80039
